$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
  2 = 1
  3 = 2
  4 = 1
  5 = 1
  6 = 2
  7 = 1
  8 = 1
  9 = 1
  10 = 1
  11 = 2
  12 = 0
  13 = 2
  14 = 1
  15 = 2
  16 = 2
  17 = 2
  18 = 0
  19 = 0
  20 = 0
  21 = 3
  22 = 2
  23 = 3
  24 = 2
  25 = 2
  26 = 1
  27 = 0
  28 = 1
  29 = 1
  30 = 1
  31 = 1
  32 = 0
  33 = 1
  34 = 0
  35 = 1
  36 = 2
  37 = 1
  38 = 1
  39 = 2
  40 = 0
  41 = 0
  42 = 3
  43 = 1
  44 = 2
  45 = 2
  46 = 0
  47 = 1
  48 = 1
  49 = 2
  50 = 2
  51 = 2
  52 = 2
  53 = 2
  54 = 1
  55 = 0
  56 = 1
  57 = 1
  58 = 2
  59 = 0
  60 = 0
  61 = 2
  62 = 2
  63 = 0
  64 = 3
  65 = 2
  67 = 3
  68 = 2
  70 = 2
  72 = 1
}

foreach ($row in $values.Keys) {
  $ws.Range("G$row").Value = $values[$row]
}
